$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37; existing rows 37:64 shift down to 38:65,
# carrying their formatting (incl. the date number format on column D).
$ws.Rows("37:37").Insert()

# Populate the newly inserted row 37 with the new weekly price-report entry.
$ws.Range("A37").Value = 10
$ws.Range("B37").Value = 'Vega Modelo de Temuco'
$ws.Range("C37").Value = 'La Araucanía'
$ws.Range("D37").Value = 44589
$ws.Range("E37").Value = 9
$ws.Range("F37").Value = 100112030
$ws.Range("G37").Value = 'Poroto granado'
$ws.Range("H37").Value = 'Sin especificar'
$ws.Range("I37").Value = 'Primera'
$ws.Range("J37").Value = 140
$ws.Range("K37").Value = 25000
$ws.Range("L37").Value = 28000
$ws.Range("M37").Value = 26179
$ws.Range("N37").Value = '$/saco 25 kilos'
$ws.Range("O37").Value = 'Región de La Araucanía'
$ws.Range("P37").Value = 1047
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = 'Hortaliza'
